$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update giang vien (lecturer) names - simulate "DeTai updated when GiangVien is updated"
$ws.Range("B6").Value = "Trần Văn Eii"
$ws.Range("B2").Value = "Trần Văn Aii."

# Update the active selection to B2, matching the diff
$ws.Range("B2").Select()
